$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new boundary-value rows (0 < x < 200) ---
$ws.Range("A5").Value = "0,4,5"
$ws.Range("B5").Value = "Not a Triangle"

$ws.Range("A6").Value = "200,4,5"
$ws.Range("B6").Value = "Not a Triangle"

# --- Set explicit column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.4
$ws.Columns.Item(2).ColumnWidth = 13.5

# --- Center (horizontal + vertical) all the used cells in one style write ---
# Build the combined alignment on a scratch cell, then copy the *format only*
# onto the target range so only a single new cell style is created.
$scratch = $ws.Range("D1")
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("A1:B6").PasteSpecial(-4122)
$scratch.Clear()

$ws.Range("C6").Select()
